$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A44").Value = "Raffaele Prezzi"
$ws.Range("B44").Value = "Elia Barozzi | I Magnifici"
$ws.Range("C44").Value = "Filippo Benetti | I Magnifici"
$ws.Range("D44").Value = "Carlo Stedile | MAI UNA GIOIA"
$ws.Range("E44").Value = "Leonardo Viola | Shark Attack"
$ws.Range("F44").Value = "Andreas Galli | SdrumALA"
